# Commit: "se hacen cambios y se agrega generador de llaves a admin"
#
# The underlying OOXML diff for this commit touches xl/workbook.xml:
#   - fileVersion/@rupBuild, x15ac:absPath/@url, xr:revisionPtr/@documentId
#     and bookViews/workbookView/@xWindow,yWindow,windowWidth,windowHeight
#     are all Excel-session/UI chrome that the authoring machine stamped
#     on save (local build number, local file path, co-authoring session
#     id, last on-screen window rect) - not workbook content, and not
#     something a script running against the workbook would set.
#   - The one real, content-level change is the sheet being renamed from
#     "22-09-2022" to "Hoja1".
#
# Apply that rename via the Excel object model.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Name = "Hoja1"
